$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.320.88'
$ws.Range('E2').Value = '  +2.38%  '
$ws.Range('D3').Value = '1.822.68'
$ws.Range('E3').Value = '  +1.61%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.77'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4678'
$ws.Range('E7').Value = '  +4.97%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3771'
$ws.Range('E8').Value = '  +2.86%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07429'
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8727'
$ws.Range('E10').Value = '  +1.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.67'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '1.820.52'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.676'
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.414'
$ws.Range('E14').Value = '  +2.83%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '92.68'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.07099'
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008797'
$ws.Range('E18').Value = '  +1.61%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.97'
$ws.Range('E20').Value = '  +1.25%  '
$ws.Range('D21').Value = '27.322.97'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.311'
$ws.Range('E22').Value = '  +3.55%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.94'
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('D24').Value = '2.045.26'
$ws.Range('E24').Value = '  -4.87%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.931'
$ws.Range('E25').Value = '  -2.75%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '151.59'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.256'
$ws.Range('E27').Value = '  +4.40%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.57'
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.310'
$ws.Range('E29').Value = '  +2.44%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '117.16'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08936'
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7814'
$ws.Range('E32').Value = '  +5.59%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.186'
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.529'
$ws.Range('E34').Value = '  +2.11%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.947'
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  +1.44%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01971'
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05244'
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5341'
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.236'
$ws.Range('E41').Value = '  +3.88%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.895'
$ws.Range('E42').Value = '  +2.35%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.359'
$ws.Range('E43').Value = '  +21.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1695'
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.621'
$ws.Range('E45').Value = '  +2.80%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5091'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.57'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '105.48'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.676'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06356'
$ws.Range('E51').Value = '  +1.03%  '
